$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.232.44'
$ws.Range("E2").Value = '  -2.22%  '
$ws.Range("D3").Value = '3.887.28'
$ws.Range("E3").Value = '  -2.32%  '
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E4").Value = '  +0.02%  '
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '597.87'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +1.83%  '
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '169.72'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +7.73%  '
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.675'
$c.Style = "Normal"
$ws.Range("E7").Value = '  -1.29%  '
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '1.00'
$c.Style = "Normal"
$ws.Range("E8").Value = '  +0.17%  '
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.756'
$c.Style = "Normal"
$ws.Range("E9").Value = '  +0.81%  '
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.177'
$c.Style = "Normal"
$ws.Range("E10").Value = '  +4.95%  '
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '53.96'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +0.31%  '
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '0.0000324'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +1.22%  '
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '11.47'
$c.Style = "Normal"
$ws.Range("E13").Value = '  +6.00%  '
$ws.Range("D14").Value = '4.512.82'
$ws.Range("E14").Value = '  -2.20%  '
$ws.Range("D15").Value = '3.890.39'
$ws.Range("E15").Value = '  -2.15%  '
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '21.00'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +3.06%  '
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '13.97'
$c.Style = "Normal"
$ws.Range("E17").Value = '  -0.64%  '
$ws.Range("E18").Value = '  -4.72%  '
$ws.Range("E19").Value = '  -1.95%  '
$ws.Range("D20").Value = '71.148.94'
$ws.Range("E20").Value = '  -1.88%  '
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '438.59'
$c.Style = "Normal"
$ws.Range("E21").Value = '  +1.12%  '
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '4.76'
$c.Style = "Normal"
$ws.Range("E22").Value = '  +1.67%  '
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '95.01'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -1.23%  '
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '3.30'
$c.Style = "Normal"
$ws.Range("E24").Value = '  -3.70%  '
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '13.92'
$c.Style = "Normal"
$ws.Range("E25").Value = '  -2.99%  '
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '4.11'
$c.Style = "Normal"
$ws.Range("E26").Value = '  -7.45%  '
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '11.33'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +2.32%  '
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '5.94'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +0.07%  '
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '10.40'
$c.Style = "Normal"
$ws.Range("E29").Value = '  -2.41%  '
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '35.29'
$c.Style = "Normal"
$ws.Range("E30").Value = '  -3.22%  '
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '8.14'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +4.15%  '
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '13.70'
$c.Style = "Normal"
$ws.Range("E32").Value = '  +0.64%  '
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '49.63'
$c.Style = "Normal"
$ws.Range("E33").Value = '  +2.38%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '0.126'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -3.92%  '
$ws.Range("B35").Value = 'OKB'
$ws.Range("C35").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '70.32'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +2.07%  '
$ws.Range("D36").Value = '0.0₃0995'
$ws.Range("E36").Value = '  +13.33%  '
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '636.60'
$c.Style = "Normal"
$ws.Range("E37").Value = '  -6.42%  '
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.426'
$c.Style = "Normal"
$ws.Range("E38").Value = '  -2.50%  '
$ws.Range("B39").Value = 'dogwifhat'
$ws.Range("C39").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '3.37'
$c.Style = "Normal"
$ws.Range("E39").Value = '  +28.88%  '
$ws.Range("B40").Value = 'Kaspa'
$ws.Range("C40").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '0.146'
$c.Style = "Normal"
$ws.Range("E40").Value = '  -0.41%  '
$ws.Range("B41").Value = 'Dai'
$ws.Range("C41").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E41").Value = '  -0.08%  '
$ws.Range("E42").Value = '  -0.04%  '
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.28'
$c.Style = "Normal"
$ws.Range("E43").Value = '  -1.33%  '
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '0.0475'
$c.Style = "Normal"
$ws.Range("E44").Value = '  -2.67%  '
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '10.18'
$c.Style = "Normal"
$ws.Range("E45").Value = '  -5.66%  '
$ws.Range("E46").Value = '  +3.48%  '
$ws.Range("E47").Value = '  -3.35%  '
$ws.Range("E48").Value = '  -15.43%  '
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '3.30'
$c.Style = "Normal"
$ws.Range("E49").Value = '  -3.15%  '
$ws.Range("D50").Value = '2.836.53'
$ws.Range("E50").Value = '  +1.73%  '
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '0.000275'
$c.Style = "Normal"
$ws.Range("E51").Value = '  +1.39%  '
